$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC (95 cell edits) ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 59499.5
$ws.Range("J3").Value = 60000
$ws.Range("L3").Value = 60000
$ws.Range("N3").Value = -60228
$ws.Range("H6").Value = 45786.09
$ws.Range("I6").Value = 50364.5
$ws.Range("K6").Value = 151093.5
$ws.Range("M6").Value = -150981.5
$ws.Range("H32").Value = 7331.391
$ws.Range("I32").Value = 5842.846
$ws.Range("K32").Value = 5842.846
$ws.Range("M32").Value = -5516.846
$ws.Range("H33").Value = 698.93335
$ws.Range("I33").Value = 622.53845
$ws.Range("K33").Value = 622.53845
$ws.Range("M33").Value = -393.53845
$ws.Range("H52").Value = 999
$ws.Range("I52").Value = 999
$ws.Range("K52").Value = 2997
$ws.Range("M52").Value = -2837
$ws.Range("H55").Value = 378.73685
$ws.Range("J55").Value = 434.2
$ws.Range("L55").Value = 434.2
$ws.Range("N55").Value = -862.2
$ws.Range("H62").Value = 3767.2307
$ws.Range("I62").Value = 3821.875
$ws.Range("J62").Value = 3679.8
$ws.Range("K62").Value = 3821.875
$ws.Range("L62").Value = 3679.8
$ws.Range("M62").Value = -3197.875
$ws.Range("N62").Value = -4927.8
$ws.Range("H65").Value = 3767.2307
$ws.Range("I65").Value = 3821.875
$ws.Range("J65").Value = 3679.8
$ws.Range("K65").Value = 19109.375
$ws.Range("L65").Value = 18399
$ws.Range("M65").Value = -15989.375
$ws.Range("N65").Value = -24639
$ws.Range("H100").Value = 2301.1052
$ws.Range("I100").Value = 2355.7334
$ws.Range("J100").Value = 2096.25
$ws.Range("K100").Value = 2355.7334
$ws.Range("L100").Value = 2096.25
$ws.Range("M100").Value = -1814.7334
$ws.Range("N100").Value = -3178.25
$ws.Range("H101").Value = 4212.125
$ws.Range("I101").Value = 1794.2222
$ws.Range("J101").Value = 7320.857
$ws.Range("K101").Value = 5382.6666
$ws.Range("L101").Value = 21962.571
$ws.Range("M101").Value = -3760.6666
$ws.Range("N101").Value = -25206.571
$ws.Range("H102").Value = 59499.5
$ws.Range("J102").Value = 60000
$ws.Range("L102").Value = 60000
$ws.Range("N102").Value = -66490
$ws.Range("H103").Value = 3325.5
$ws.Range("I103").Value = 4001
$ws.Range("K103").Value = 12003
$ws.Range("M103").Value = -11417
$ws.Range("H113").Value = 111152296
$ws.Range("I113").Value = 250002130
$ws.Range("K113").Value = 250002130
$ws.Range("M113").Value = -249998876
$ws.Range("H121").Value = 1429.2333
$ws.Range("J121").Value = 1449.2858
$ws.Range("L121").Value = 4347.857400000001
$ws.Range("N121").Value = -7841.857400000001
$ws.Range("H125").Value = 1603486.1
$ws.Range("I125").Value = 2565182.5
$ws.Range("K125").Value = 23086642.5
$ws.Range("M125").Value = -23084182.5
$ws.Range("H132").Value = 1809.2174
$ws.Range("I132").Value = 1736.909
$ws.Range("K132").Value = 5210.727000000001
$ws.Range("M132").Value = -2680.727000000001
$ws.Range("H134").Value = 109988.5
$ws.Range("J134").Value = 109988.5
$ws.Range("L134").Value = 109988.5
$ws.Range("N134").Value = -120128.5
$ws.Range("H135").Value = 1606.4193
$ws.Range("I135").Value = 1105.5927
$ws.Range("K135").Value = 9950.334299999999
$ws.Range("M135").Value = -7415.334299999999
$ws.Range("H136").Value = 183331.67
$ws.Range("J136").Value = 183331.67
$ws.Range("L136").Value = 183331.67
$ws.Range("N136").Value = -193531.67
$ws.Range("H141").Value = 1795
$ws.Range("I141").Value = 1837.5
$ws.Range("J141").Value = 1200
$ws.Range("K141").Value = 5512.5
$ws.Range("L141").Value = 3600
$ws.Range("M141").Value = -332.5
$ws.Range("N141").Value = -13960

# ---- Sheet: ARM (86 cell edits) ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1648
$ws.Range("I2").Value = 1216.25
$ws.Range("J2").Value = 2991.2222
$ws.Range("K2").Value = 1216.25
$ws.Range("L2").Value = 2991.2222
$ws.Range("M2").Value = -1103.25
$ws.Range("N2").Value = -3217.2222
$ws.Range("H26").Value = 11258
$ws.Range("I26").Value = 11258
$ws.Range("K26").Value = 11258
$ws.Range("M26").Value = -10928
$ws.Range("H43").Value = 35171.582
$ws.Range("J43").Value = 33469.3
$ws.Range("L43").Value = 33469.3
$ws.Range("N43").Value = -34095.3
$ws.Range("H45").Value = 76937.22
$ws.Range("I45").Value = 101960.35
$ws.Range("J45").Value = 5442.5713
$ws.Range("K45").Value = 101960.35
$ws.Range("L45").Value = 5442.5713
$ws.Range("M45").Value = -101583.35
$ws.Range("N45").Value = -6196.5713
$ws.Range("H63").Value = 3542.1428
$ws.Range("I63").Value = 3542.1428
$ws.Range("K63").Value = 3542.1428
$ws.Range("M63").Value = -2856.1428
$ws.Range("H66").Value = 3542.1428
$ws.Range("I66").Value = 3542.1428
$ws.Range("K66").Value = 17710.714
$ws.Range("M66").Value = -14278.714
$ws.Range("H74").Value = 811.1111
$ws.Range("I74").Value = 811.1111
$ws.Range("K74").Value = 811.1111
$ws.Range("M74").Value = 62.88890000000004
$ws.Range("H77").Value = 811.1111
$ws.Range("I77").Value = 811.1111
$ws.Range("K77").Value = 4055.5555
$ws.Range("M77").Value = 312.4445000000001
$ws.Range("H97").Value = 930.4828
$ws.Range("I97").Value = 964.8
$ws.Range("J97").Value = 716
$ws.Range("K97").Value = 964.8
$ws.Range("L97").Value = 716
$ws.Range("M97").Value = -468.8
$ws.Range("N97").Value = -1708
$ws.Range("H102").Value = 3643.037
$ws.Range("I102").Value = 1744.2273
$ws.Range("K102").Value = 1744.2273
$ws.Range("M102").Value = -122.2273
$ws.Range("H110").Value = 2415.1155
$ws.Range("I110").Value = 520.75
$ws.Range("J110").Value = 8729.666999999999
$ws.Range("K110").Value = 520.75
$ws.Range("L110").Value = 8729.666999999999
$ws.Range("M110").Value = 1524.25
$ws.Range("N110").Value = -12819.667
$ws.Range("H116").Value = 1648
$ws.Range("I116").Value = 1216.25
$ws.Range("J116").Value = 2991.2222
$ws.Range("K116").Value = 1216.25
$ws.Range("L116").Value = 2991.2222
$ws.Range("M116").Value = 1077.75
$ws.Range("N116").Value = -7579.2222
$ws.Range("H122").Value = 2110.4707
$ws.Range("I122").Value = 2074.2307
$ws.Range("K122").Value = 6222.6921
$ws.Range("M122").Value = -3772.6921
$ws.Range("H124").Value = 36124.75
$ws.Range("J124").Value = 36124.75
$ws.Range("L124").Value = 36124.75
$ws.Range("N124").Value = -45944.75
$ws.Range("H125").Value = 60000
$ws.Range("J125").Value = 60000
$ws.Range("L125").Value = 60000
$ws.Range("N125").Value = -69840
$ws.Range("H132").Value = 2902.4614
$ws.Range("I132").Value = 2961.0435
$ws.Range("J132").Value = 2453.3333
$ws.Range("K132").Value = 8883.130500000001
$ws.Range("L132").Value = 7359.999899999999
$ws.Range("M132").Value = -6353.130500000001
$ws.Range("N132").Value = -12419.9999
$ws.Range("H133").Value = 79999
$ws.Range("J133").Value = 79999
$ws.Range("L133").Value = 79999
$ws.Range("N133").Value = -85059

# ---- Sheet: BSM (40 cell edits) ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1648
$ws.Range("I3").Value = 1216.25
$ws.Range("J3").Value = 2991.2222
$ws.Range("K3").Value = 1216.25
$ws.Range("L3").Value = 2991.2222
$ws.Range("M3").Value = -1102.25
$ws.Range("N3").Value = -3219.2222
$ws.Range("H64").Value = 458.77777
$ws.Range("I64").Value = 446
$ws.Range("J64").Value = 465.16666
$ws.Range("K64").Value = 446
$ws.Range("L64").Value = 465.16666
$ws.Range("M64").Value = -221
$ws.Range("N64").Value = -915.16666
$ws.Range("H67").Value = 458.77777
$ws.Range("I67").Value = 446
$ws.Range("J67").Value = 465.16666
$ws.Range("K67").Value = 446
$ws.Range("L67").Value = 465.16666
$ws.Range("M67").Value = 334
$ws.Range("N67").Value = -2025.16666
$ws.Range("H99").Value = 1557.7273
$ws.Range("I99").Value = 1017.125
$ws.Range("K99").Value = 1017.125
$ws.Range("M99").Value = 480.875
$ws.Range("H105").Value = 3100
$ws.Range("I105").Value = 2150
$ws.Range("K105").Value = 2150
$ws.Range("M105").Value = -403
$ws.Range("H107").Value = 1798.5
$ws.Range("I107").Value = 1134.6364
$ws.Range("J107").Value = 4232.6665
$ws.Range("K107").Value = 1134.6364
$ws.Range("L107").Value = 4232.6665
$ws.Range("M107").Value = 785.3635999999999
$ws.Range("N107").Value = -8072.6665
$ws.Range("H134").Value = 3465
$ws.Range("I134").Value = 2939.2083
$ws.Range("K134").Value = 8817.624899999999
$ws.Range("M134").Value = -6282.624899999999

# ---- Sheet: CRP (65 cell edits) ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5236.75
$ws.Range("I31").Value = 3240.6667
$ws.Range("J31").Value = 11225
$ws.Range("K31").Value = 3240.6667
$ws.Range("L31").Value = 11225
$ws.Range("M31").Value = -2945.6667
$ws.Range("N31").Value = -11815
$ws.Range("H34").Value = 5236.75
$ws.Range("I34").Value = 3240.6667
$ws.Range("J34").Value = 11225
$ws.Range("K34").Value = 3240.6667
$ws.Range("L34").Value = 11225
$ws.Range("M34").Value = -3038.6667
$ws.Range("N34").Value = -11629
$ws.Range("H50").Value = 22000
$ws.Range("J50").Value = 34000
$ws.Range("L50").Value = 34000
$ws.Range("N50").Value = -35250
$ws.Range("H51").Value = 19425.8
$ws.Range("J51").Value = 33995
$ws.Range("L51").Value = 33995
$ws.Range("N51").Value = -35467
$ws.Range("H58").Value = 7957.0835
$ws.Range("I58").Value = 4855.4287
$ws.Range("K58").Value = 4855.4287
$ws.Range("M58").Value = -4652.4287
$ws.Range("H59").Value = 43330
$ws.Range("J59").Value = 44995
$ws.Range("L59").Value = 44995
$ws.Range("N59").Value = -47285
$ws.Range("H60").Value = 43666.332
$ws.Range("J60").Value = 32000
$ws.Range("L60").Value = 32000
$ws.Range("N60").Value = -33022
$ws.Range("H61").Value = 19425.8
$ws.Range("J61").Value = 33995
$ws.Range("L61").Value = 33995
$ws.Range("N61").Value = -34691
$ws.Range("H62").Value = 5469.25
$ws.Range("I62").Value = 4730.8
$ws.Range("K62").Value = 4730.8
$ws.Range("M62").Value = -4106.8
$ws.Range("H65").Value = 5469.25
$ws.Range("I65").Value = 4730.8
$ws.Range("K65").Value = 23654
$ws.Range("M65").Value = -20534
$ws.Range("H122").Value = 3144.7083
$ws.Range("I122").Value = 3218.3333
$ws.Range("J122").Value = 3022
$ws.Range("K122").Value = 9654.999899999999
$ws.Range("L122").Value = 9066
$ws.Range("M122").Value = -7204.999899999999
$ws.Range("N122").Value = -13966
$ws.Range("H132").Value = 387967.47
$ws.Range("I132").Value = 590173.8
$ws.Range("K132").Value = 1770521.4
$ws.Range("M132").Value = -1767991.4
$ws.Range("H136").Value = 7957.0835
$ws.Range("I136").Value = 4855.4287
$ws.Range("K136").Value = 14566.2861
$ws.Range("M136").Value = -12016.2861
$ws.Range("H141").Value = 296484.12
$ws.Range("J141").Value = 358105.16
$ws.Range("L141").Value = 358105.16
$ws.Range("N141").Value = -368465.16

# ---- Sheet: CUL (23 cell edits) ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 6508605
$ws.Range("I44").Value = 13015000
$ws.Range("K44").Value = 39045000
$ws.Range("M44").Value = -39044602
$ws.Range("H47").Value = 100
$ws.Range("I47").Value = 100
$ws.Range("K47").Value = 300
$ws.Range("M47").Value = 131
$ws.Range("H107").Value = 2322.2727
$ws.Range("I107").Value = 1001
$ws.Range("K107").Value = 3003
$ws.Range("M107").Value = -1083
$ws.Range("H113").Value = 2503.5557
$ws.Range("I113").Value = 824.5
$ws.Range("K113").Value = 2473.5
$ws.Range("M113").Value = -303.5
$ws.Range("H122").Value = 5289.857
$ws.Range("I122").Value = 4092.75
$ws.Range("J122").Value = 5644.5557
$ws.Range("K122").Value = 36834.75
$ws.Range("L122").Value = 50801.0013
$ws.Range("M122").Value = -34384.75
$ws.Range("N122").Value = -55701.0013

# ---- Sheet: GSM (56 cell edits) ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 11661.777
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 11661.777
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 11661.777
$ws.Range("M70").Value = $null
$ws.Range("N70").Value = -12201.777
$ws.Range("H73").Value = 11661.777
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 11661.777
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 11661.777
$ws.Range("M73").Value = $null
$ws.Range("N73").Value = -13533.777
$ws.Range("H80").Value = 6708.5
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").Value = $null
$ws.Range("H83").Value = 6708.5
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").Value = $null
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").Value = $null
$ws.Range("H97").Value = 1023.913
$ws.Range("I97").Value = 999
$ws.Range("J97").Value = 1113.6
$ws.Range("K97").Value = 999
$ws.Range("L97").Value = 1113.6
$ws.Range("M97").Value = -503
$ws.Range("N97").Value = -2105.6
$ws.Range("H102").Value = 5075.8
$ws.Range("I102").Value = 4645.96
$ws.Range("K102").Value = 4645.96
$ws.Range("M102").Value = -3023.96
$ws.Range("H113").Value = 3542.8572
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").Value = $null
$ws.Range("H122").Value = 2524.76
$ws.Range("I122").Value = 2657.1052
$ws.Range("J122").Value = 2105.6667
$ws.Range("K122").Value = 7971.3156
$ws.Range("L122").Value = 6317.000100000001
$ws.Range("M122").Value = -5521.3156
$ws.Range("N122").Value = -11217.0001
$ws.Range("H132").Value = 1787.3462
$ws.Range("I132").Value = 1871.1305
$ws.Range("K132").Value = 5613.3915
$ws.Range("M132").Value = -3083.3915
$ws.Range("H141").Value = 268792.44
$ws.Range("J141").Value = 268792.44
$ws.Range("L141").Value = 268792.44
$ws.Range("N141").Value = -279152.44

# ---- Sheet: LTW (12 cell edits) ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").Value = $null
$ws.Range("H40").Value = 8112.125
$ws.Range("I40").Value = 8299
$ws.Range("K40").Value = 8299
$ws.Range("M40").Value = -8163
$ws.Range("H93").Value = 11360
$ws.Range("I93").Value = 2162.7334
$ws.Range("K93").Value = 2162.7334
$ws.Range("M93").Value = -914.7334000000001

# ---- Sheet: WVR (37 cell edits) ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 24333
$ws.Range("J24").Value = 24333
$ws.Range("L24").Value = 24333
$ws.Range("N24").Value = -24793
$ws.Range("H100").Value = 669.3333
$ws.Range("I100").Value = 332
$ws.Range("K100").Value = 664
$ws.Range("M100").Value = -123
$ws.Range("H122").Value = 5062.1763
$ws.Range("I122").Value = 3382.875
$ws.Range("J122").Value = 6554.8887
$ws.Range("K122").Value = 10148.625
$ws.Range("L122").Value = 19664.6661
$ws.Range("M122").Value = -7698.625
$ws.Range("N122").Value = -24564.6661
$ws.Range("H132").Value = 140020.33
$ws.Range("I132").Value = 184585.05
$ws.Range("J132").Value = 3850.3333
$ws.Range("K132").Value = 553755.1499999999
$ws.Range("L132").Value = 11550.9999
$ws.Range("M132").Value = -551225.1499999999
$ws.Range("N132").Value = -16610.9999
$ws.Range("H136").Value = 2612.319
$ws.Range("I136").Value = 1738.9474
$ws.Range("K136").Value = 5216.8422
$ws.Range("M136").Value = -2666.8422
$ws.Range("H139").Value = 80122.5
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 80122.5
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 80122.5
$ws.Range("M139").Value = $null
$ws.Range("N139").Value = -90402.5
$ws.Range("H140").Value = 65823.57000000001
$ws.Range("J140").Value = 65823.57000000001
$ws.Range("L140").Value = 65823.57000000001
$ws.Range("N140").Value = -76183.57000000001
